$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Locate the paragraph that currently reads "... (901.28, see
#    901.47, 901.73)" inside the "Pre-flight Planning" table row.
# -----------------------------------------------------------------
$searchText = "901.28, see 901.47, 901.73"
$paraCount = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*$searchText*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    Write-Output "Target paragraph not found - aborting"
} else {
    # -----------------------------------------------------------------
    # 2. Trim the citation down to just "(901.28)" - the "see 901.47,
    #    901.73" part moves to its own new bullet about aerodrome
    #    distance.
    # -----------------------------------------------------------------
    $targetPara = $d.Paragraphs.Item($targetIndex)
    $findRange = $targetPara.Range
    $replaced = $findRange.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "901.28", 2)
    Write-Output ("Replaced citation text: " + $replaced)

    # -----------------------------------------------------------------
    # 3. Split off a new paragraph right after it and fill it with the
    #    new "Consider the distance from any aerodrome ..." sentence.
    # -----------------------------------------------------------------
    $targetPara = $d.Paragraphs.Item($targetIndex)
    $targetPara.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Range.Text = "Consider the distance from any aerodrome of the intended operation (901.47, 901.73)"
    Write-Output "Inserted new aerodrome-distance paragraph"

    # -----------------------------------------------------------------
    # 4. The extra line bumped the row heights of the "Pre-flight
    #    Planning" row and (via reflow) the following "Site Survey"
    #    row - set them to their new values explicitly.
    # -----------------------------------------------------------------
    $table = $d.Tables.Item(1)
    for ($r = 1; $r -le $table.Rows.Count; $r++) {
        $row = $table.Rows.Item($r)
        $label = $row.Cells.Item(1).Range.Text
        if ($label -like "*Pre-flight Planning*") {
            $row.Height = 4751 / 20
            Write-Output ("Set Pre-flight Planning row height to " + $row.Height)
        } elseif ($label -like "*Site Survey*") {
            $row.Height = 4238 / 20
            Write-Output ("Set Site Survey row height to " + $row.Height)
        }
    }
}
